$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1431
$ws1.Range("F3").Value = 3014
$ws1.Range("F4").Value = 34
$ws1.Range("F5").Value = 265
$ws1.Range("F6").Value = 282

# Sheet "全部类型" (4th sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1431
$ws4.Range("F3").Value = 3014
$ws4.Range("F4").Value = 34
$ws4.Range("F5").Value = 265
$ws4.Range("F7").Value = 282
